$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are written as plain text (avoid numeric
# auto-coercion e.g. trailing-zero stripping like "7.10" -> "7.1").

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.294.89"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.294.11"
$ws.Range("E3").Value = "  +1.39%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.64"
$ws.Range("E5").Value = "  +1.18%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.21"
$ws.Range("E6").Value = "  -1.01%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.291.18"
$ws.Range("E8").Value = "  +1.48%  "

# Row 9
$ws.Range("E9").Value = "  -0.69%  "

# Row 10
$ws.Range("E10").Value = "  +0.72%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.48"
$ws.Range("E11").Value = "  +1.81%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.464"
$ws.Range("E12").Value = "  -0.14%  "

# Row 13
$ws.Range("E13").Value = "  -1.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.28"
$ws.Range("E14").Value = "  +0.17%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.836.43"
$ws.Range("E15").Value = "  +1.37%  "

# Row 16
$ws.Range("E16").Value = "  +1.27%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.295.44"
$ws.Range("E17").Value = "  +1.40%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.384.64"
$ws.Range("E18").Value = "  +0.01%  "

# Row 19
$ws.Range("E19").Value = "  +0.53%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.75"
$ws.Range("E20").Value = "  +0.18%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.88"
$ws.Range("E21").Value = "  -2.07%  "

# Row 22
$ws.Range("E22").Value = "  +0.23%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.92"
$ws.Range("E23").Value = "  -0.46%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.79"
$ws.Range("E24").Value = "  +4.69%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.76"
$ws.Range("E25").Value = "  +1.30%  "

# Row 26
$ws.Range("E26").Value = "  -0.03%  "

# Row 27
$ws.Range("E27").Value = "  +0.65%  "

# Row 28
$ws.Range("E28").Value = "  +0.04%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.10"
$ws.Range("E29").Value = "  -1.90%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.06"
$ws.Range("E30").Value = "  -0.26%  "

# Row 31
$ws.Range("E31").Value = "  -0.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.39"
$ws.Range("E32").Value = "  +2.86%  "

# Row 33
$ws.Range("E33").Value = "  -2.43%  "

# Row 34
$ws.Range("E34").Value = "  -1.43%  "

# Row 35
$ws.Range("E35").Value = "  -0.10%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("E36").Value = "  +0.89%  "

# Row 37
$ws.Range("E37").Value = "  -0.85%  "

# Row 38
$ws.Range("E38").Value = "  +2.17%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0396"
$ws.Range("E39").Value = "  +0.98%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.109.93"
$ws.Range("E40").Value = "  +4.20%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "426.03"
$ws.Range("E41").Value = "  +1.39%  "

# Row 42
$ws.Range("E42").Value = "  +7.20%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.24"
$ws.Range("E43").Value = "  -1.28%  "

# Row 44
$ws.Range("E44").Value = "  -2.38%  "

# Row 45
$ws.Range("E45").Value = "  -1.78%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.18"
$ws.Range("E46").Value = "  +0.62%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "128.11"
$ws.Range("E47").Value = "  +4.31%  "

# Row 48
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.03%  "

# Row 49
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.87"
$ws.Range("E49").Value = "  +7.78%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.04"
$ws.Range("E50").Value = "  +0.96%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.113"
$ws.Range("E51").Value = "  -1.12%  "
